$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric strings need to be forced to Text
# so Excel does not auto-convert them to numbers (matching the source inlineStr type).
$textForceRefs = @("D5","D6","D7","D9","D10","D11","D12","D15","D21","D22","D23","D25","D26","D27","D28","D29","D30","D31","D34","D35","D36","D40","D42","D43","D44","D45","D46","D47","D48","D49","D51")
foreach ($r in $textForceRefs) { $ws.Range($r).NumberFormat = "@" }

$ws.Range("D2").Value = "43.978.14"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "2.262.99"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "0.661"
$ws.Range("E5").Value = "  +6.02%  "
$ws.Range("D6").Value = "233.87"
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("D7").Value = "63.56"
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.454"
$ws.Range("E9").Value = "  +7.99%  "
$ws.Range("D10").Value = "0.0982"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").Value = "57.99"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("D12").Value = "26.91"
$ws.Range("E12").Value = "  +5.34%  "
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("D14").Value = "2.601.91"
$ws.Range("D15").Value = "15.66"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("E16").Value = "  +5.79%  "
$ws.Range("E17").Value = "  +4.01%  "
$ws.Range("D18").Value = "2.266.97"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "43.917.37"
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("D20").Value = "0.0₃0987"
$ws.Range("E20").Value = "  +2.54%  "
$ws.Range("D21").Value = "74.19"
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("D22").Value = "6.17"
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("D23").Value = "251.04"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("D25").Value = "2.46"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").Value = "3.33"
$ws.Range("E26").Value = "  +19.32%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "2.23"
$ws.Range("E27").Value = "  -3.70%  "
$ws.Range("D28").Value = "9.95"
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("D29").Value = "22.28"
$ws.Range("E29").Value = "  +9.16%  "
$ws.Range("D30").Value = "174.06"
$ws.Range("E30").Value = "  +1.92%  "
$ws.Range("D31").Value = "0.137"
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("E33").Value = "  +5.05%  "
$ws.Range("D34").Value = "5.03"
$ws.Range("E34").Value = "  +8.18%  "
$ws.Range("D35").Value = "0.0686"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D36").Value = "5.00"
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("E37").Value = "  -2.14%  "
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").Value = "0.0255"
$ws.Range("E40").Value = "  +4.35%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "8.82"
$ws.Range("E42").Value = "  +6.06%  "
$ws.Range("D43").Value = "0.000224"
$ws.Range("E43").Value = "  +6.92%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "98.85"
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "17.27"
$ws.Range("E45").Value = "  +3.38%  "
$ws.Range("D46").Value = "0.0954"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").Value = "1.19"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").Value = "4.38"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").Value = "2.34"
$ws.Range("E49").Value = "  +1.98%  "
$ws.Range("D50").Value = "1.455.00"
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("D51").Value = "10.00"
$ws.Range("E51").Value = "  -3.10%  "

foreach ($r in $textForceRefs) { $ws.Range($r).Style = "Normal" }
